$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.182548999786377
$ws.Range("B1").Value = 1.251620054244995
$ws.Range("C1").Value = 1.416848421096802
$ws.Range("D1").Value = 2.251186847686768
$ws.Range("E1").Value = 4.486071586608887
